$wb = $excel.ActiveWorkbook

# --- Overview sheet: b.md row status -> "Ready for handoff" for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) gets a new handoff file + datetime, status updated ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-09 12:27:29"

# --- de-de sheet: b.md row (row 3) gets a new handoff file + datetime, status updated ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-09 12:27:34"
